$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PutUsers")

# New "Address" column header + values (order matters for shared-string table ordering)
$ws.Cells.Item(1, 8).Value = "Address"
$ws.Cells.Item(2, 8).Value = '{"Address1":"12234","Address2":"ABC Apt 12","City":"Farmington Hills","State":"MI","Country":"USA"}'
$ws.Cells.Item(3, 8).Value = '{"Address1":"Main Rd","Address2":"Roswell","City":"Ketty","State":"TX","Country":"USA"}'
$ws.Cells.Item(4, 8).Value = '{"Address1":"Main Rd","Address2":"Roswell","City":"Farmington Hills","State":"MI","Country":"USA"}'
$ws.Cells.Item(5, 8).Value = '{"Address1":"Main Rd","Address2":"Apt 12","City":"Atlanta","State":"GA","Country":"USA"}'

# Other cell edits
$ws.Cells.Item(3, 5).Value = "Veg"
$ws.Cells.Item(6, 2).Value = "Potterr"
$ws.Cells.Item(5, 2).Value = "Dek"
$ws.Cells.Item(5, 6).Value = "Nuts"

# Remaining address value for row 6
$ws.Cells.Item(6, 8).Value = '{"Address1":"Main","Address2":"Apt 12","City":"Tampa","State":"FL","Country":"USA"}'

# Update selection to match target
$ws.Range("G8").Select()
